# Applies the "Added requirements for input config." change to the
# MASTER SPREADSHEET worksheet of the requirements workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Rows 271-277 previously had an empty "Status" (column E) cell.
#    They are now marked "Incomplete", matching the rest of the table.
# ---------------------------------------------------------------------
$incompleteRows = 271,272,273,274,275,276,277
foreach ($r in $incompleteRows) {
    $ws.Cells.Item($r, 5).Value = "Incomplete"
}

# ---------------------------------------------------------------------
# 2. Rows 278-283 were blank placeholder rows; they now hold six new
#    "SYS-INP-*" input-configuration requirements.
# ---------------------------------------------------------------------
$newRows = @(
    @{ Row=278; Id="SYS-INP-001";       Priority="HIGH"; Text="The program shall read the input elements and configure the input source automatically."; RestyleE=$false },
    @{ Row=279; Id="SYS-INP-002";       Priority="HIGH"; Text="The program shall allow input sources from SOCKET, QTJSBSIM, or NONE interfaces.";           RestyleE=$true  },
    @{ Row=280; Id="SYS-INP-003";       Priority="LOW";  Text="The program shall allow the data rate for the input source to manage timing of data collection."; RestyleE=$true },
    @{ Row=281; Id="SYS-INP-004-001";   Priority="HIGH"; Text="The program shall specify a network port through a network port number";                     RestyleE=$true  },
    @{ Row=282; Id="SYS-INP-004-002";   Priority="HIGH"; Text="The network port number MUST be a POSITIVE INTEGER.";                                        RestyleE=$true  },
    @{ Row=283; Id="SYS-INP-005-001";   Priority="HIGH"; Text="The program shall have an ""Action"" attribute to the input souce.";                          RestyleE=$true  }
)

foreach ($item in $newRows) {
    $r = $item.Row

    $ws.Cells.Item($r, 2).Value = $item.Id
    $ws.Cells.Item($r, 4).Value = $item.Text

    # Column C (priority) needs a style change (fill colour) in addition
    # to its value, so copy the formatting of an existing HIGH/LOW cell
    # before writing the value.
    if ($item.Priority -eq "HIGH") {
        $ws.Range("C271").Copy()
    } else {
        $ws.Range("C276").Copy()
    }
    $ws.Range($ws.Cells.Item($r, 3), $ws.Cells.Item($r, 3)).PasteSpecial(-4122)
    $ws.Cells.Item($r, 3).Value = $item.Priority

    # Column E (status): row 278 already carries the "Incomplete" fill,
    # the rest need it copied over before the value is written.
    if ($item.RestyleE) {
        $ws.Range("E271").Copy()
        $ws.Range($ws.Cells.Item($r, 5), $ws.Cells.Item($r, 5)).PasteSpecial(-4122)
    }
    $ws.Cells.Item($r, 5).Value = "Incomplete"
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. The author left the selection on the last newly entered cell.
# ---------------------------------------------------------------------
$ws.Range("D283").Select()
